$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-13 Saturday" "2025-09-14 Sunday"

Replace-Text "57×14=" "84×43="
Replace-Text "34×37=" "75×77="
Replace-Text "26×33=" "69×86="
Replace-Text "51×88=" "74×34="
Replace-Text "63×71=" "40×27="
Replace-Text "11×68=" "64×41="
Replace-Text "31×21=" "37×40="
Replace-Text "73×24=" "16×28="
Replace-Text "33×66=" "41×85="
Replace-Text "92×74=" "48×63="
Replace-Text "35×57=" "23×19="
Replace-Text "43×37=" "79×86="
Replace-Text "79×64=" "13×20="
Replace-Text "87×46=" "11×66="
Replace-Text "20×90=" "57×12="
Replace-Text "35×96=" "61×69="
Replace-Text "73×13=" "51×13="
Replace-Text "92×50=" "19×30="
Replace-Text "38×29=" "45×17="
Replace-Text "19×76=" "38×78="
Replace-Text "49×71=" "59×66="
Replace-Text "20×45=" "21×71="
Replace-Text "27×12=" "21×87="
Replace-Text "49×91=" "11×18="
Replace-Text "84×27=" "54×50="
